# Horarios actualizados Linea 141 - 1009
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with a
# new scrape timestamp and refreshed rows of arrival data.

$wb = $excel.ActiveWorkbook

$scrapTime = "04:21:09"

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: $scrapTime"
$ws1.Cells.Item(3,1).Value = "Total filas: 10"

$rows1 = @(
    @($scrapTime, "04:46", "215A_EL PATO", 25, "LP1912"),
    @($scrapTime, "04:53", "11_ETCHEVERRY", 32, "LP1912"),
    @($scrapTime, "05:16", "17_ROMERO", 55, "LP1912"),
    @($scrapTime, "05:22", "23_HERNANDEZ", 61, "LP1912"),
    @($scrapTime, "05:35", "215B_EL PATO", 74, "LP1912"),
    @($scrapTime, "05:46", "15_ABASTO", 85, "LP1912"),
    @($scrapTime, "05:54", "10_OLMOS", 93, "LP1912"),
    @($scrapTime, "06:04", "16_SANTA ANA", 103, "LP1912"),
    @($scrapTime, "06:11", "215A_EL PATO", 110, "LP1912"),
    @($scrapTime, "06:14", "225_HARAS DEL SUR", 113, "LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r,1).Value = $row[0]
    $ws1.Cells.Item($r,2).Value = $row[1]
    $ws1.Cells.Item($r,3).Value = $row[2]
    $ws1.Cells.Item($r,4).Value = $row[3]
    $ws1.Cells.Item($r,5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: $scrapTime"
$ws2.Cells.Item(3,1).Value = "Total filas: 3"

$ws2.Cells.Item(5,1).Value = "Hora_Scrap"
$ws2.Cells.Item(5,2).Value = "Hora_Llegada"
$ws2.Cells.Item(5,3).Value = "Linea"
$ws2.Cells.Item(5,4).Value = "Minutos"
$ws2.Cells.Item(5,5).Value = "Parada"

$rows2 = @(
    @($scrapTime, "04:46", "215A_EL PATO", 25, "LP1912"),
    @($scrapTime, "05:35", "215B_EL PATO", 74, "LP1912"),
    @($scrapTime, "06:11", "215A_EL PATO", 110, "LP1912")
)

$r = 6
foreach ($row in $rows2) {
    $ws2.Cells.Item($r,1).Value = $row[0]
    $ws2.Cells.Item($r,2).Value = $row[1]
    $ws2.Cells.Item($r,3).Value = $row[2]
    $ws2.Cells.Item($r,4).Value = $row[3]
    $ws2.Cells.Item($r,5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: $scrapTime"
$ws3.Cells.Item(3,1).Value = "Total filas: 2"

$ws3.Cells.Item(5,1).Value = "Hora_Scrap"
$ws3.Cells.Item(5,2).Value = "Hora_Llegada"
$ws3.Cells.Item(5,3).Value = "Linea"
$ws3.Cells.Item(5,4).Value = "Minutos"
$ws3.Cells.Item(5,5).Value = "Parada"

$rows3 = @(
    @($scrapTime, "05:44", "215A_LA PLATA", 83, "L6173"),
    @($scrapTime, "06:09", "215A_LA PLATA", 108, "L6173")
)

$r = 6
foreach ($row in $rows3) {
    $ws3.Cells.Item($r,1).Value = $row[0]
    $ws3.Cells.Item($r,2).Value = $row[1]
    $ws3.Cells.Item($r,3).Value = $row[2]
    $ws3.Cells.Item($r,4).Value = $row[3]
    $ws3.Cells.Item($r,5).Value = $row[4]
    $r++
}
